$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59: paper_chipping
$ws.Range("A59").Value = "paper_chipping"
$ws.Range("B59").Value = "paper"
$ws.Range("C59").Value = "Wood Handling"
$ws.Range("D59").Value = "chips"
$ws.Range("E59").Value = "outflow"
$ws.Range("F59").Value = "data/paper/paper_var.xlsx"
$ws.Range("G59").Value = "Wood"
$ws.Range("H59").Value = "data/paper/paper_calc.xlsx"
$ws.Range("I59").Value = "Wood"

# Row 60: paper_pulping
$ws.Range("A60").Value = "paper_pulping"
$ws.Range("B60").Value = "paper"
$ws.Range("C60").Value = "Pulp Preparation"
$ws.Range("D60").Value = "unbleached pulp"
$ws.Range("E60").Value = "outflow"
$ws.Range("F60").Value = "data/paper/paper_var.xlsx"
$ws.Range("G60").Value = "Pulp"
$ws.Range("H60").Value = "data/paper/paper_calc.xlsx"
$ws.Range("I60").Value = "Pulp"

# Apply the "Text" number format (matches style used across column A,C,E,F,G,H,I in existing rows)
$ws.Range("A59:A60").NumberFormat = "@"
$ws.Range("C59:C60").NumberFormat = "@"
$ws.Range("E59:E60").NumberFormat = "@"
$ws.Range("F59:F60").NumberFormat = "@"
$ws.Range("G59:G60").NumberFormat = "@"
$ws.Range("H59:H60").NumberFormat = "@"
$ws.Range("I59:I60").NumberFormat = "@"

# D59 uses the smaller (11pt) font style seen elsewhere in column D; D60 keeps default style
$ws.Range("D59").Font.Size = 11

# Update selection to match the new active cell shown in the saved file
$ws.Range("F62").Select()
